# "corrected fullrun to spikein"
# The "purpose" column (G) on Sheet1 was mislabeled "fullRNASeq" for every
# data row; replace it with "spikein" for all 52 data rows (rows 2-53).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = 53
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = "spikein"
}

# Match the author's final selection in the saved file (G3:G53, active cell G3).
$ws.Range("G3:G53").Select()
